# Delete rows that are no longer present in the updated public exposure sites list.
# Deleting from the bottom up keeps the remaining row numbers valid as we go.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Delete()   # Wonthaggi
$ws.Rows.Item(24).Delete()   # Mordialloc
$ws.Rows.Item(23).Delete()   # Moorabbin
$ws.Rows.Item(9).Delete()    # Hallam
$ws.Rows.Item(3).Delete()    # Brighton
